$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Replace embedded CRLF line breaks in the descriptive labels (column E)
#     with literal "<br>" tags, matching the updated plotting HTML. ---
$ws.Range("E2:E13").Value2  = "**Global**:<br>Implemented by<br>All other countries"
$ws.Range("E14:E25").Value2 = "**High-income**:<br>All other HICs and<br>not some MICs (such as China)"
$ws.Range("E26:E37").Value2 = "**International**:<br>Some countries (e.g. EU, UK, Brazil)<br>and not others (e.g. U.S., China)"

# --- Updated mean / CI_low / CI_high estimates (refreshed ggplot2 3.5.1 run) ---
$ws.Range("B2").Value2 = 73.6710361597934
$ws.Range("C2").Value2 = 72.3033300855294
$ws.Range("D2").Value2 = 75.0387422340573

$ws.Range("B12").Value2 = 77.5532714938846
$ws.Range("C12").Value2 = 73.1395709765103
$ws.Range("D12").Value2 = 81.966972011259

$ws.Range("B14").Value2 = 68.7235111211832
$ws.Range("C14").Value2 = 67.2867144309077
$ws.Range("D14").Value2 = 70.1603078114586

$ws.Range("B24").Value2 = 69.8291634341244
$ws.Range("C24").Value2 = 64.9837035538317
$ws.Range("D24").Value2 = 74.6746233144171
